$wb = $excel.ActiveWorkbook

# --- Update the Date on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-11-16T21:33:27+00:00"

# --- Fix the swapped MFTH/MMTH codes on the "Concepts" sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("B24").Value = "MMTH"
$concepts.Range("B25").Value = "MFTH"
